{"js": "const body = context.document.body;\nconst pairs = [\n  [\"283\u00d75=1415\", \"481\u00d79=4329\"],\n  [\"669\u00d72=1338\", \"436\u00d72=872\"],\n  [\"433\u00d78=3464\", \"665\u00d76=3990\"],\n  [\"617\u00d74=2468\", \"315\u00d77=2205\"],\n  [\"576\u00d79=5184\", \"296\u00d75=1480\"],\n  [\"636\u00d73=1908\", \"858\u00d73=2574\"],\n  [\"787\u00d75=3935\", \"940\u00d73=2820\"],\n  [\"898\u00d74=3592\", \"481\u00d78=3848\"],\n  [\"807\u00d74=3228\", \"716\u00d75=3580\"],\n  [\"511\u00d78=4088\", \"704\u00d79=6336\"],\n  [\"232\u00d72=464\", \"649\u00d74=2596\"],\n  [\"689\u00d78=5512\", \"946\u00d77=6622\"],\n  [\"562\u00d75=2810\", \"301\u00d79=2709\"],\n  [\"523\u00d73=1569\", \"695\u00d72=1390\"],\n  [\"292\u00d77=2044\", \"194\u00d75=970\"],\n  [\"495\u00d74=1980\", \"506\u00d72=1012\"],\n  [\"698\u00d73=2094\", \"588\u00d73=1764\"],\n  [\"774\u00d76=4644\", \"140\u00d73=420\"],\n  [\"128\u00d74=512\", \"405\u00d79=3645\"],\n  [\"950\u00d79=8550\", \"108\u00d77=756\"],\n  [\"544\u00d76=3264\", \"986\u00d72=1972\"],\n  [\"293\u00d73=879\", \"240\u00d76=1440\"],\n  [\"389\u00d78=3112\", \"562\u00d73=1686\"],\n  [\"298\u00d79=2682\", \"481\u00d77=3367\"],\n  [\"656\u00d73=1968\", \"716\u00d78=5728\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    ,@(\"283\u00d75=1415\", \"481\u00d79=4329\")\n    ,@(\"669\u00d72=1338\", \"436\u00d72=872\")\n    ,@(\"433\u00d78=3464\", \"665\u00d76=3990\")\n    ,@(\"617\u00d74=2468\", \"315\u00d77=2205\")\n    ,@(\"576\u00d79=5184\", \"296\u00d75=1480\")\n    ,@(\"636\u00d73=1908\", \"858\u00d73=2574\")\n    ,@(\"787\u00d75=3935\", \"940\u00d73=2820\")\n    ,@(\"898\u00d74=3592\", \"481\u00d78=3848\")\n    ,@(\"807\u00d74=3228\", \"716\u00d75=3580\")\n    ,@(\"511\u00d78=4088\", \"704\u00d79=6336\")\n    ,@(\"232\u00d72=464\", \"649\u00d74=2596\")\n    ,@(\"689\u00d78=5512\", \"946\u00d77=6622\")\n    ,@(\"562\u00d75=2810\", \"301\u00d79=2709\")\n    ,@(\"523\u00d73=1569\", \"695\u00d72=1390\")\n    ,@(\"292\u00d77=2044\", \"194\u00d75=970\")\n    ,@(\"495\u00d74=1980\", \"506\u00d72=1012\")\n    ,@(\"698\u00d73=2094\", \"588\u00d73=1764\")\n    ,@(\"774\u00d76=4644\", \"140\u00d73=420\")\n    ,@(\"128\u00d74=512\", \"405\u00d79=3645\")\n    ,@(\"950\u00d79=8550\", \"108\u00d77=756\")\n    ,@(\"544\u00d76=3264\", \"986\u00d72=1972\")\n    ,@(\"293\u00d73=879\", \"240\u00d76=1440\")\n    ,@(\"389\u00d78=3112\", \"562\u00d73=1686\")\n    ,@(\"298\u00d79=2682\", \"481\u00d77=3367\")\n    ,@(\"656\u00d73=1968\", \"716\u00d78=5728\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
